$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPLKKPS005-001")
$ws.Activate()

# A new "SIDEBAR_SUBMENU_SUBMENU" level is introduced in the sidebar
# navigation hierarchy. Insert a new column at M (13) so everything from
# the old KODE_PARAMETER column onward shifts one column to the right.
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).ColumnWidth = 14.17

# Populate the new column's header and data first (so the default style
# Excel assigns to freshly written values doesn't get clobbered by the
# format copy below).
$ws.Range("M1").Value = "SIDEBAR_SUBMENU_SUBMENU"
$ws.Range("M2").Value = "Setup Jenis Parameter"

# Column L (SIDEBAR_SUBMENU) now holds the higher-level submenu label,
# while the previous value ("Setup Jenis Parameter") moved into M2 above.
$ws.Range("L2").Value = "Setup Kelengkapan Kepesertaan"

# Match the formatting of the neighbouring cells: M2 picks up N2's plain
# style, L2 picks up K2's style (quote-prefixed left/center alignment).
$ws.Range("N2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# The GRUP_PARAMETER value for row 2 (now column R, after the insert) is
# removed entirely.
$ws.Range("R2").Clear()

# Refresh the window selection to match the new layout.
$ws.Range("S2").Select() | Out-Null
